$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet "SearchCases". Add the new
# "AddNewCases" sheet right after it - this becomes the active sheet/tab,
# matching activeTab moving from 0 to 1 and SearchCases' tabSelected
# flipping to false.
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "AddNewCases"

# Header row
$newSheet.Range("A1").Value = "Cpimsid"
$newSheet.Range("B1").Value = "FirstName"
$newSheet.Range("C1").Value = "LastName"
$newSheet.Range("D1").Value = "AssesmentDueDate"
$newSheet.Range("E1").Value = "age"

# Data row
$newSheet.Range("A2").Value = "Id_cpims_1"
$newSheet.Range("B2").Value = "Janani"
$newSheet.Range("C2").Value = "Panchalingam"

# Keep "10-10-2024" as literal text (not auto-converted to a date serial):
# temporarily force a text number format while assigning the value, then
# restore General so no stray formatting is left behind.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "10-10-2024"
$newSheet.Range("D2").NumberFormat = "General"

$newSheet.Range("E2").Value = 33
